$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.424.71"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.238.90"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +1.52%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.93"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.20"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.570"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.65"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0802"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.19"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.304.32"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.00%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.55"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "44.107.80"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.33"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.50%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.82%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "237.22"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.98"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.01%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.36%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.32%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.93"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "153.24"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0793"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.05%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.04%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.111"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.88%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "14.93"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.86%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0298"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.59%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.792.18"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.14%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "78.61"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -8.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "70.16"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "98.43"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.48%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.07"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.39"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.22%  "
